# Update worksheet values to reflect the new TPM-derived computation.
# Only the numeric columns E..H and M..T for data rows 2-5 change;
# columns A-D, I-L remain identical (identifiers / counts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{
        E = 2
        F = 0.6666666666666666
        G = 0.2438943333333334
        H = 0.7316830000000001
        M = 5.273684
        N = 15.821052
        O = 0.0510821201937383
        P = 0.0510821201937383
        Q = 1.286221643390667
        R = 11.575994790516
        S = 0.0510821201937383
        T = 0.0510821201937383
    }
    3 = @{
        E = 2
        F = 0.6666666666666666
        G = 0.2438943333333334
        H = 0.7316830000000001
        O = 0.5598845502029881
        P = 0.5598845502029881
        Q = 14.09760643332511
        R = 126.878457899926
        S = 0.5598845502029881
        T = 0.5598845502029881
    }
    4 = @{
        E = 2
        F = 0.6666666666666666
        G = 0.2438943333333334
        H = 0.7316830000000001
        M = 32.95839133333334
        N = 98.87517400000002
        O = 0.3192425840231603
        P = 0.3192425840231604
        Q = 8.038364881982446
        R = 72.34528393784203
        S = 0.3192425840231603
        T = 0.3192425840231604
    }
    5 = @{
        E = 2
        F = 0.6666666666666666
        G = 0.2438943333333334
        H = 0.7316830000000001
        M = 7.205150000000001
        N = 21.61545
        O = 0.06979074558011317
        P = 0.06979074558011318
        Q = 1.757295255816667
        R = 15.81565730235
        S = 0.06979074558011317
        T = 0.06979074558011318
    }
}

foreach ($rowNum in $data.Keys) {
    $rowData = $data[$rowNum]
    foreach ($col in $rowData.Keys) {
        $addr = "$col$rowNum"
        $ws.Range($addr).Value = $rowData[$col]
    }
}

$wb.Save()
